$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "position" -> "response"
$ws.Range("B1").Value = "response"

# Row 2 keeps cue_file value (images/YellowIcon.png), response becomes Incorrect
$ws.Range("A2").Value = "images/YellowIcon.png"
$ws.Range("B2").Value = "Incorrect"

# Row 3 becomes OrangeIcon / Incorrect
$ws.Range("A3").Value = "images/OrangeIcon.png"
$ws.Range("B3").Value = "Incorrect"

# Row 4 becomes StopIcon / Correct
$ws.Range("A4").Value = "images/StopIcon.png"
$ws.Range("B4").Value = "Correct"

# Remove the now-unused trailing practice rows (previously rows 5-7)
$ws.Rows("5:7").Delete()

[void]$ws.Range("E24").Select()
